$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1598
$ws1.Range("F3").Value = 259
$ws1.Range("F4").Value = 8506
$ws1.Range("F6").Value = 61
$ws1.Range("F7").Value = 19
$ws1.Range("F9").Value = 1323
$ws1.Range("F10").Value = 92
$ws1.Range("F13").Value = 9156
$ws1.Range("F14").Value = 149
$ws1.Range("F16").Value = 209
$ws1.Range("F17").Value = 168
$ws1.Range("F18").Value = 339
$ws1.Range("F19").Value = 6014
$ws1.Range("F20").Value = 1038
$ws1.Range("F21").Value = 56
$ws1.Range("F22").Value = 37
$ws1.Range("F23").Value = 95

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1598
$ws4.Range("F4").Value = 259
$ws4.Range("F5").Value = 8506
$ws4.Range("F7").Value = 61
$ws4.Range("F8").Value = 19
$ws4.Range("F10").Value = 1323
$ws4.Range("F11").Value = 92
$ws4.Range("F16").Value = 9156
$ws4.Range("F17").Value = 149
$ws4.Range("F19").Value = 209
$ws4.Range("F20").Value = 168
$ws4.Range("F21").Value = 339
$ws4.Range("F22").Value = 6014
$ws4.Range("F23").Value = 1038
$ws4.Range("F24").Value = 56
$ws4.Range("F25").Value = 37
$ws4.Range("F26").Value = 95
